$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.230.04"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").Value = "1.858.90"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'0.7024"
$ws.Range("E5").Value = "  +0.06%  "

$ws.Range("D6").Value = "'237.46"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").Value = "'0.9993"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'0.07939"
$ws.Range("E8").Value = "  +7.31%  "

$ws.Range("D9").Value = "'0.3050"
$ws.Range("E9").Value = "  +0.68%  "

$ws.Range("D10").Value = "'23.30"
$ws.Range("E10").Value = "  -0.14%  "

$ws.Range("D11").Value = "'0.08193"
$ws.Range("E11").Value = "  +0.99%  "

$ws.Range("D12").Value = "1.845.36"
$ws.Range("E12").Value = "  -0.70%  "

$ws.Range("D13").Value = "'0.7192"
$ws.Range("E13").Value = "  -0.81%  "

$ws.Range("D14").Value = "'5.172"
$ws.Range("E14").Value = "  -0.66%  "

$ws.Range("D15").Value = "'89.13"
$ws.Range("E15").Value = "  +0.10%  "

$ws.Range("D16").Value = "29.234.36"
$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("D17").Value = "'5.777"
$ws.Range("E17").Value = "  +0.03%  "

$ws.Range("D18").Value = "'13.36"
$ws.Range("E18").Value = "  +2.83%  "

$ws.Range("D19").Value = "'0.000007778"
$ws.Range("E19").Value = "  +1.67%  "

$ws.Range("D20").Value = "'236.77"
$ws.Range("E20").Value = "  -1.99%  "

$ws.Range("D21").Value = "'0.9997"
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").Value = "2.106.61"
$ws.Range("E22").Value = "  +0.75%  "

$ws.Range("D23").Value = "'0.9991"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").Value = "'7.456"
$ws.Range("E24").Value = "  -1.69%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'9.004"
$ws.Range("E25").Value = "  +0.49%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'161.80"
$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.1465"
$ws.Range("E27").Value = "  -0.49%  "

$ws.Range("D28").Value = "'18.04"
$ws.Range("E28").Value = "  +0.14%  "

$ws.Range("E29").Value = "  +4.03%  "

$ws.Range("D30").Value = "'1.434"
$ws.Range("E30").Value = "  +4.21%  "

$ws.Range("D31").Value = "'4.423"
$ws.Range("E31").Value = "  -0.66%  "

$ws.Range("D32").Value = "'1.484"
$ws.Range("E32").Value = "  -0.23%  "

$ws.Range("D33").Value = "'4.051"
$ws.Range("E33").Value = "  +1.27%  "

$ws.Range("D34").Value = "'0.05221"
$ws.Range("E34").Value = "  +0.40%  "

$ws.Range("D35").Value = "'1.171"
$ws.Range("E35").Value = "  -1.01%  "

$ws.Range("D36").Value = "'0.7083"
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("D38").Value = "'2.671"
$ws.Range("E38").Value = "  +0.98%  "

$ws.Range("D39").Value = "'0.01847"
$ws.Range("E39").Value = "  -1.13%  "

$ws.Range("D40").Value = "'2.718"
$ws.Range("E40").Value = "  +1.75%  "

$ws.Range("D41").Value = "'0.9268"
$ws.Range("E41").Value = "  +2.47%  "

$ws.Range("D42").Value = "1.138.16"
$ws.Range("E42").Value = "  +8.39%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4279"
$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.922"
$ws.Range("E44").Value = "  +0.54%  "

$ws.Range("D45").Value = "'70.80"
$ws.Range("E45").Value = "  +0.99%  "

$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("D47").Value = "'103.45"
$ws.Range("E47").Value = "  +1.93%  "

$ws.Range("D48").Value = "'1.797"
$ws.Range("E48").Value = "  +2.60%  "

$ws.Range("D49").Value = "2.003.02"
$ws.Range("E49").Value = "  +0.39%  "

$ws.Range("D50").Value = "'9.190"
$ws.Range("E50").Value = "  +0.15%  "

$ws.Range("D51").Value = "'6.997"
$ws.Range("E51").Value = "  -1.67%  "
